$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Update the ROOT version number shown on the schema slide (Rectangle 181):
# text changes from "5.34.34" to "5.34.36" and the run is made bold.
$shape = $s.Shapes.Item(118)
$tr = $shape.TextFrame.TextRange
$tr.Text = "5.34.36"
$tr.Font.Bold = $true
